$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 75
$ws.Range("F2").Value = 56
$ws.Range("H2").Value = 56
$ws.Range("E3").Value = 28
$ws.Range("E4").Value = 32
$ws.Range("F4").Value = 21
$ws.Range("H4").Value = 21
$ws.Range("E5").Value = 95
$ws.Range("F5").Value = 59
$ws.Range("H5").Value = 59
$ws.Range("E8").Value = 6
$ws.Range("E10").Value = 319
$ws.Range("F10").Value = 153
$ws.Range("H10").Value = 153
$ws.Range("E11").Value = 218
$ws.Range("E12").Value = 322
$ws.Range("F14").Value = 42
$ws.Range("H14").Value = 42
$ws.Range("E16").Value = 136
$ws.Range("E17").Value = 59
$ws.Range("E19").Value = 13
$ws.Range("E21").Value = 105
$ws.Range("E22").Value = 121
$ws.Range("E23").Value = 127
$ws.Range("E24").Value = 141
$ws.Range("E25").Value = 166
$ws.Range("F25").Value = 80
$ws.Range("H25").Value = 80
$ws.Range("E26").Value = 90
$ws.Range("E27").Value = 213
$ws.Range("F27").Value = 109
$ws.Range("H27").Value = 109
$ws.Range("E28").Value = 129
$ws.Range("F28").Value = 46
$ws.Range("H28").Value = 46
$ws.Range("E29").Value = 127
$ws.Range("F29").Value = 75
$ws.Range("H29").Value = 75
$ws.Range("E30").Value = 149
$ws.Range("F30").Value = 83
$ws.Range("H30").Value = 83
$ws.Range("E32").Value = 131
$ws.Range("E33").Value = 204
$ws.Range("E34").Value = 149
$ws.Range("E35").Value = 98
$ws.Range("F35").Value = 61
$ws.Range("H35").Value = 61
$ws.Range("E37").Value = 109
$ws.Range("E38").Value = 64
$ws.Range("E39").Value = 134
$ws.Range("E40").Value = 183
$ws.Range("F40").Value = 85
$ws.Range("H40").Value = 85
$ws.Range("E41").Value = 263
$ws.Range("F41").Value = 119
$ws.Range("H41").Value = 119
$ws.Range("E42").Value = 243
$ws.Range("F42").Value = 127
$ws.Range("H42").Value = 127
$ws.Range("E43").Value = 78
$ws.Range("E44").Value = 216
$ws.Range("F44").Value = 112
$ws.Range("H44").Value = 112
$ws.Range("E45").Value = 88
$ws.Range("E46").Value = 207
$ws.Range("F46").Value = 114
$ws.Range("H46").Value = 114
$ws.Range("E47").Value = 304
$ws.Range("F47").Value = 153
$ws.Range("H47").Value = 153
$ws.Range("E48").Value = 140
$ws.Range("E49").Value = 186
$ws.Range("E50").Value = 159
$ws.Range("E51").Value = 150
